$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (was "SCD0337", now "SCD0025")
$ws.Name = "SCD0025"

# Update the TC_ID cell from "DGS-352" to "SCD0025-007"
$ws.Range("B2").Value = "SCD0025-007"

# The TC_ID column now holds a wider value, so its best-fit width grows
# (was 9 chars, now ~12.4 chars wide to fit "SCD0025-007")
$ws.Columns.Item(2).ColumnWidth = 11.66

# Move the selection cursor to B3 (matches the saved selection state)
$ws.Range("B3").Select()
